# GW_Project_Planning.xlsx - "Reorg files and folders"
# Re-work the task list: extend/rewrite several comments, insert new task
# rows (co2/temperature chart planning, folder clean-up), and add a new
# "Questions for Peleke on 11/30" section at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert three new rows after row 2 (becomes rows 3-5) and one more
#    after the (now shifted) row 11 (becomes row 12). Excel copies the
#    formatting of the row above automatically.
# ---------------------------------------------------------------------
$ws.Rows("3:5").Insert()
$ws.Rows("12").Insert()

# ---------------------------------------------------------------------
# 2. Row 1 - header row (text unchanged, only the shared-string identity
#    of "Hours" changes internally - no user-visible change needed).
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Hours"

# ---------------------------------------------------------------------
# 3. Row 2 - extend the comment, bump the row height 45 -> 60.
# ---------------------------------------------------------------------
$ws.Range("C2").Value = "Need update from Aurora and probably meeting with Peleke on 11/27." + [char]10 + "Look at making a co2 chart similar to  "
$ws.Rows(2).RowHeight = 60

# ---------------------------------------------------------------------
# 4. New rows 3-5 (co2 / temperature planning)
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "Clean up folders"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Chris and Aurora need to organize and clean up folders and files."
$ws.Rows(3).RowHeight = 30

$ws.Range("A4").Value = "co2 chart siimilar to area plot for final presentation."
$ws.Range("B4").Value = "N/A"
$ws.Range("C4").Value = "Aurora to choose countries for co2 and temperature and send to chris."
$ws.Rows(4).RowHeight = 30

$ws.Range("A5").Value = "temperature chart to show multiple temperature by country as a final chart."
$ws.Range("B5").Value = "N/A"
$ws.Range("C5").Value = "Look for an example chart for temperature."
$ws.Rows(5).RowHeight = 30

# ---------------------------------------------------------------------
# 5. Rows 6-11 (former rows 3-8) - update the comment text in a couple of
#    spots and bump row 7's height.
# ---------------------------------------------------------------------
$ws.Range("C6").Value = 'This is the README_Project-GW.md file that is displayed on a page when the "About" link on the home page is clicked.'

$ws.Range("C7").Value = 'This is the file that sits in the root directory of Project-GW GitHub repo, and it will be displayed to anyone visiting the repository.  This could have the same information as the "About" readme file or could contain additional details.'
$ws.Rows(7).RowHeight = 105

# ---------------------------------------------------------------------
# 6. New row 12 (Matplotlib / D3 animated chart - nice to have)
# ---------------------------------------------------------------------
$ws.Range("A12").Value = "Use Matplotlib widgets, D3+, or D3 to create an animated chart."
$ws.Range("B12").Value = "N/A"
$ws.Range("C12").Value = "This is a nice to have."
$ws.Rows(12).RowHeight = 30

# ---------------------------------------------------------------------
# 7. New "Questions for Peleke on 11/30:" section (rows 16-20), formatted
#    like the header row (yellow fill, medium bottom border, wrap text).
# ---------------------------------------------------------------------
$ws.Range("A16").Value = "Questions for Peleke on 11/30:"
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = ""

$hdr = $ws.Range("A16:C16")
$hdr.Interior.Color = 65535
$hdr.Borders.Item(9).LineStyle = 1
$hdr.Borders.Item(9).Weight = -4138
$hdr.WrapText = $true
$ws.Range("A16").Font.Bold = $true
$ws.Rows(16).RowHeight = 15.75

$ws.Range("A17").Value = "Choropleth world and USA maps files availability?"
$ws.Range("C17").Value = "Possible world map with pop-up data like level1 d3 homework."
$ws.Rows(17).RowHeight = 30

$ws.Range("A18").Value = "Temperature chart example similar to co2 example?"
$ws.Rows(18).RowHeight = 30

$ws.Range("A19").Value = "Histograms - what are they supposed to show us?  What parameters to use in creating them?"
$ws.Range("C19").Value = "For future learning and possible inclusion in a page of how we got there."
$ws.Rows(19).RowHeight = 45

$ws.Range("A20").Value = "Box plots - what are they supposed to show us?  What parameters to use in creating them?"
$ws.Range("C20").Value = "For future learning and possible inclusion in a page of how we got there."
$ws.Rows(20).RowHeight = 45

# ---------------------------------------------------------------------
# 8. Column widths / selection.
# ---------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 47.5703125
$ws.Columns("B").ColumnWidth = 11.42578125

$ws.Range("A2").Select()

Write-Output "edit complete"
